# Update dredge model-selection table: replace "Cruise" predictor with
# "Month" predictor, shifting column roles (B=Depth, C=DRM, D=Month,
# E=Depth:DRM, F=Depth:Month, G=DRM:Month) and refreshing the per-row
# coefficients/statistics that came out of the refit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Depth"
$ws.Range("C1").Value = "DRM"
$ws.Range("D1").Value = "Month"
$ws.Range("E1").Value = "Depth:DRM"
$ws.Range("F1").Value = "Depth:Month"
$ws.Range("G1").Value = "DRM:Month"
$ws.Range("A2").Value = 12.91990777828679
$ws.Range("B2").Value = -2.046248793187126
$ws.Range("C2").Value = -2.176953027767409
$ws.Range("D2").Value = "+"
$ws.Range("E2").Value = $null
$ws.Range("F2").Value = "+"
$ws.Range("B3").Value = -1.515338749643475
$ws.Range("C3").Value = -1.343390044136112
$ws.Range("D3").Value = "+"
$ws.Range("E3").Value = 1.847617620814608
$ws.Range("F3").Value = "+"
$ws.Range("G3").Value = $null
$ws.Range("L3").Value = 1.106901876489331
$ws.Range("M3").Value = 0.1909227713804208
$ws.Range("N3").Value = 0.585417580799876
$ws.Range("A4").Value = 13.1530185801983
$ws.Range("B4").Value = -4.37180474184414
$ws.Range("C4").Value = -2.833740903033822
$ws.Range("D4").Value = "+"
$ws.Range("B5").Value = -1.097475592281952
$ws.Range("C5").Value = $null
$ws.Range("D5").Value = "+"
$ws.Range("E5").Value = $null
$ws.Range("F5").Value = "+"
$ws.Range("M5").Value = 0.08072389957981176
$ws.Range("B6").Value = -2.368967872894364
$ws.Range("C6").Value = -2.917429547357891
$ws.Range("D6").Value = "+"
$ws.Range("E6").Value = $null
$ws.Range("G6").Value = "+"
$ws.Range("M6").Value = 0.08006847493661134
$ws.Range("B7").Value = -4.689977246746887
$ws.Range("C7").Value = -4.000685560660525
$ws.Range("D7").Value = "+"
$ws.Range("F7").Value = $null
$ws.Range("G7").Value = "+"
$ws.Range("K7").Value = 211.0185507675862
$ws.Range("L7").Value = 3.808198152255954
$ws.Range("B8").Value = -4.349391754018686
$ws.Range("C8").Value = -2.52827661896059
$ws.Range("D8").Value = "+"
$ws.Range("E8").Value = 0.808466702645554
$ws.Range("G8").Value = $null
$ws.Range("B9").Value = -1.740688236563396
$ws.Range("C9").Value = -1.841999608292806
$ws.Range("D9").Value = "+"
$ws.Range("E9").Value = 1.758970135312606
$ws.Range("G9").Value = "+"
$ws.Range("H9").Value = 0.6518858270538
$ws.Range("M9").Value = 0.03557937496615352
$ws.Range("B10").Value = -4.624341065898574
$ws.Range("C10").Value = -2.726464850910515
$ws.Range("D10").Value = $null
$ws.Range("M10").Value = 0.02204094873307984
$ws.Range("B11").Value = -4.495683064867013
$ws.Range("C11").Value = -1.989216486691722
$ws.Range("D11").Value = $null
$ws.Range("E11").Value = 2.032559186829854
$ws.Range("G11").Value = $null
$ws.Range("M11").Value = 0.01951368851393099
$ws.Range("B12").Value = -4.654095223407779
$ws.Range("C12").Value = -3.680889026665232
$ws.Range("D12").Value = "+"
$ws.Range("E12").Value = 0.6814816759814242
$ws.Range("F12").Value = $null
$ws.Range("G12").Value = "+"
$ws.Range("B13").Value = -4.246929022316333
$ws.Range("C13").Value = $null
$ws.Range("D13").Value = "+"
$ws.Range("K13").Value = 215.3647379974934
$ws.Range("M13").Value = 0.005630100287908304
$ws.Range("B14").Value = -4.49021478784986
$ws.Range("C14").Value = $null
$ws.Range("M14").Value = 0.00221280419567206
$ws.Range("B15").Value = $null
$ws.Range("C15").Value = -2.642935265258226
$ws.Range("D15").Value = "+"
$ws.Range("M15").Value = 0.00005200773744122311
$ws.Range("B16").Value = $null
$ws.Range("D16").Value = "+"
$ws.Range("H16").Value = 0.1230646296830999
$ws.Range("M16").Value = 0.00001679921760234023
$ws.Range("N16").Value = 0.09477639193094189
$ws.Range("B17").Value = $null
$ws.Range("C17").Value = -1.811790456648805
$ws.Range("D17").Value = "+"
$ws.Range("F17").Value = $null
$ws.Range("G17").Value = "+"
$ws.Range("M17").Value = 0.00001521660675518298
$ws.Range("C18").Value = -2.498974056231954
$ws.Range("D18").Value = $null
$ws.Range("M18").Value = 0.00001260580899913714
$ws.Range("M19").Value = 0.000006477357874265096
